# Auto-generated Excel COM-interop script to apply the cryptos.xlsx update
# described by the commit "Updated symbol list on Sat Dec 24 22:56:55 UTC 2022 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric-looking "Price" column (D) values must stay stored as literal text ---
# (exactly like the original workbook, e.g. "244.52" rather than the float 244.52),
# so we temporarily force a text number format ("@") on each target cell while
# assigning its value, then restore the original "General" format afterwards.

$priceCells = $ws.Range("D2,D5,D6,D7,D8,D9,D10,D11,D12,D13,D14,D15,D16,D17,D18,D19,D20,D21,D22,D25,D26,D27,D40,D41,D42,D43,D44,D45,D47,D48,D50")
foreach ($area in $priceCells.Areas) {
    $area.NumberFormat = "@"
}

$ws.Range("D2").Value = '244.52'
$ws.Range("D5").Value = '0.06037'
$ws.Range("D6").Value = '3.396'
$ws.Range("D7").Value = '0.8136'
$ws.Range("D8").Value = '0.9205'
$ws.Range("D9").Value = '0.1439'
$ws.Range("D10").Value = '0.07459'
$ws.Range("D11").Value = '0.03377'
$ws.Range("D12").Value = '0.03065'
$ws.Range("D13").Value = '0.09414'
$ws.Range("D14").Value = '4.007'
$ws.Range("D15").Value = '0.001600'
$ws.Range("D16").Value = '0.04804'
$ws.Range("D17").Value = '0.0005942'
$ws.Range("D18").Value = '0.005622'
$ws.Range("D19").Value = '0.004154'
$ws.Range("D20").Value = '0.0009908'
$ws.Range("D21").Value = '3.653'
$ws.Range("D22").Value = '6.427'
$ws.Range("D25").Value = '0.1322'
$ws.Range("D26").Value = '0.00008505'
$ws.Range("D27").Value = '0.0002901'
$ws.Range("D40").Value = '0.03991'
$ws.Range("D41").Value = '0.1074'
$ws.Range("D42").Value = '0.002711'
$ws.Range("D43").Value = '0.003047'
$ws.Range("D44").Value = '0.006371'
$ws.Range("D45").Value = '0.00005240'
$ws.Range("D47").Value = '0.8102'
$ws.Range("D48").Value = '0.002319'
$ws.Range("D50").Value = '0.01010'

foreach ($area in $priceCells.Areas) {
    $area.NumberFormat = "General"
}

# --- Plain text columns (Coin name / Link / Volume label) ---

$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'MandalaExchangeToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("E10").Value = '9MandalaExchangeTokenMDX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("E12").Value = '11BitrueCoinBTR'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("E13").Value = '12BitMartTokenBMX'
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("E14").Value = '13MCDexMCB'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("E15").Value = '14BitForexTokenBF'
$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("E16").Value = '15CoinExTokenCET'
$ws.Range("B17").Value = 'One'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E17").Value = '16OneONE'
$ws.Range("E27").Value = '26UpBotsUBXTBestin24h'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'CEJI'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("E42").Value = '41CEJICEJI'
$ws.Range("B43").Value = 'KickToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("E43").Value = '42KickTokenKICK'
